$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$wb.Worksheets.Item(1).Name = "Through 2021-09-12"

# Update the row label text for September in column A (row 10)
$ws.Range("A10").Value = "September (through 09-12)"

# Update the September row (row 10) values
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 24
$ws.Range("F10").Value = 28
$ws.Range("G10").Value = 39
$ws.Range("H10").Value = 58

# Update the Total row (row 11) values
$ws.Range("B11").Value = 205
$ws.Range("C11").Value = 401
$ws.Range("D11").Value = 581
$ws.Range("E11").Value = 514
$ws.Range("F11").Value = 377
$ws.Range("G11").Value = 823
$ws.Range("H11").Value = 1129
